$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all")

# Row 9: new expense entry - folic acid, bal krishna sir
$ws.Range("A6").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = 66435
$ws.Range("B9").Value = "folic acid, bal krishna sir"
$ws.Range("C9").Formula = "=30+1000"

# Row 7: append ", bal krishna sir" to the description, and add 1000 to the formula
$ws.Range("B7").Value = "dinesh vinaju & prabin chiya, irika bus fare, Photocopy, apple, chocolate, bal krishna sir"
$ws.Range("C7").Formula = "=65+100+50+325+50+1000"

# Row 10: add date only
$ws.Range("A6").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = 66436

# Row 11: add date only
$ws.Range("A6").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = 66437

$excel.CutCopyMode = 0

# Update selection to B6
$ws.Range("B6").Select()
